$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header fix: I1 should match H1 formatting (style used to have border/numfmt, now matches H1) ---
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1").PasteSpecial(-4122) | Out-Null

# --- Column I (Railway Crossing) data cells: normalize style (drop stray numberFormat flag) ---
$iRange = $ws.Range("I2:I152")
$iRange.Style = "Normal"
$iRange.Font.Size = 12

# --- Column H (Station) dispatch-to-block numbers ---
# Cells that take the integer block-number style (numFmt "0", centered) - style index 6
$ws.Range("H3").Value = 1
$ws.Range("H3").NumberFormat = "0"
$ws.Range("H3").HorizontalAlignment = -4108
$ws.Range("H6").Value = 4
$ws.Range("H6").NumberFormat = "0"
$ws.Range("H6").HorizontalAlignment = -4108
$ws.Range("H9").Value = 7
$ws.Range("H9").NumberFormat = "0"
$ws.Range("H9").HorizontalAlignment = -4108
$ws.Range("H12").Value = 10
$ws.Range("H12").NumberFormat = "0"
$ws.Range("H12").HorizontalAlignment = -4108
$ws.Range("H14").Value = 12
$ws.Range("H14").NumberFormat = "0"
$ws.Range("H14").HorizontalAlignment = -4108
$ws.Range("H17").Value = 15
$ws.Range("H17").NumberFormat = "0"
$ws.Range("H17").HorizontalAlignment = -4108
$ws.Range("H19").Value = 17
$ws.Range("H19").NumberFormat = "0"
$ws.Range("H19").HorizontalAlignment = -4108
$ws.Range("H22").Value = 20
$ws.Range("H22").NumberFormat = "0"
$ws.Range("H22").HorizontalAlignment = -4108
$ws.Range("H27").Value = 25
$ws.Range("H27").NumberFormat = "0"
$ws.Range("H27").HorizontalAlignment = -4108
$ws.Range("H29").Value = 27
$ws.Range("H29").NumberFormat = "0"
$ws.Range("H29").HorizontalAlignment = -4108
$ws.Range("H32").Value = 30
$ws.Range("H32").NumberFormat = "0"
$ws.Range("H32").HorizontalAlignment = -4108
$ws.Range("H34").Value = 32
$ws.Range("H34").NumberFormat = "0"
$ws.Range("H34").HorizontalAlignment = -4108
$ws.Range("H37").Value = 35
$ws.Range("H37").NumberFormat = "0"
$ws.Range("H37").HorizontalAlignment = -4108
$ws.Range("H40").Value = 38
$ws.Range("H40").NumberFormat = "0"
$ws.Range("H40").HorizontalAlignment = -4108
$ws.Range("H43").Value = 41
$ws.Range("H43").NumberFormat = "0"
$ws.Range("H43").HorizontalAlignment = -4108
$ws.Range("H46").Value = 44
$ws.Range("H46").NumberFormat = "0"
$ws.Range("H46").HorizontalAlignment = -4108
$ws.Range("H49").Value = 47
$ws.Range("H49").NumberFormat = "0"
$ws.Range("H49").HorizontalAlignment = -4108
$ws.Range("H52").Value = 50
$ws.Range("H52").NumberFormat = "0"
$ws.Range("H52").HorizontalAlignment = -4108
$ws.Range("H55").Value = 53
$ws.Range("H55").NumberFormat = "0"
$ws.Range("H55").HorizontalAlignment = -4108
$ws.Range("H58").Value = 56
$ws.Range("H58").NumberFormat = "0"
$ws.Range("H58").HorizontalAlignment = -4108
$ws.Range("H61").Value = 59
$ws.Range("H61").NumberFormat = "0"
$ws.Range("H61").HorizontalAlignment = -4108
$ws.Range("H64").Value = 62
$ws.Range("H64").NumberFormat = "0"
$ws.Range("H64").HorizontalAlignment = -4108
$ws.Range("H70").Value = 68
$ws.Range("H70").NumberFormat = "0"
$ws.Range("H70").HorizontalAlignment = -4108
$ws.Range("H73").Value = 71
$ws.Range("H73").NumberFormat = "0"
$ws.Range("H73").HorizontalAlignment = -4108
$ws.Range("H76").Value = 74
$ws.Range("H76").NumberFormat = "0"
$ws.Range("H76").HorizontalAlignment = -4108
$ws.Range("H80").Value = 78
$ws.Range("H80").NumberFormat = "0"
$ws.Range("H80").HorizontalAlignment = -4108
$ws.Range("H84").Value = 82
$ws.Range("H84").NumberFormat = "0"
$ws.Range("H84").HorizontalAlignment = -4108
$ws.Range("H88").Value = 86
$ws.Range("H88").NumberFormat = "0"
$ws.Range("H88").HorizontalAlignment = -4108
$ws.Range("H92").Value = 90
$ws.Range("H92").NumberFormat = "0"
$ws.Range("H92").HorizontalAlignment = -4108
$ws.Range("H96").Value = 94
$ws.Range("H96").NumberFormat = "0"
$ws.Range("H96").HorizontalAlignment = -4108
$ws.Range("H100").Value = 98
$ws.Range("H100").NumberFormat = "0"
$ws.Range("H100").HorizontalAlignment = -4108
$ws.Range("H104").Value = 102
$ws.Range("H104").NumberFormat = "0"
$ws.Range("H104").HorizontalAlignment = -4108
$ws.Range("H108").Value = 106
$ws.Range("H108").NumberFormat = "0"
$ws.Range("H108").HorizontalAlignment = -4108
$ws.Range("H112").Value = 110
$ws.Range("H112").NumberFormat = "0"
$ws.Range("H112").HorizontalAlignment = -4108
$ws.Range("H120").Value = 118
$ws.Range("H120").NumberFormat = "0"
$ws.Range("H120").HorizontalAlignment = -4108
$ws.Range("H124").Value = 122
$ws.Range("H124").NumberFormat = "0"
$ws.Range("H124").HorizontalAlignment = -4108
$ws.Range("H128").Value = 126
$ws.Range("H128").NumberFormat = "0"
$ws.Range("H128").HorizontalAlignment = -4108
$ws.Range("H132").Value = 130
$ws.Range("H132").NumberFormat = "0"
$ws.Range("H132").HorizontalAlignment = -4108
$ws.Range("H136").Value = 134
$ws.Range("H136").NumberFormat = "0"
$ws.Range("H136").HorizontalAlignment = -4108
$ws.Range("H140").Value = 138
$ws.Range("H140").NumberFormat = "0"
$ws.Range("H140").HorizontalAlignment = -4108
$ws.Range("H144").Value = 142
$ws.Range("H144").NumberFormat = "0"
$ws.Range("H144").HorizontalAlignment = -4108
$ws.Range("H148").Value = 146
$ws.Range("H148").NumberFormat = "0"
$ws.Range("H148").HorizontalAlignment = -4108
$ws.Range("H152").Value = 150
$ws.Range("H152").NumberFormat = "0"
$ws.Range("H152").HorizontalAlignment = -4108

# Cells that change from wrap-center to no-wrap-center (previously matched station style)
$ws.Range("H21").Value = 19
$ws.Range("H21").WrapText = $false
$ws.Range("H72").Value = 70
$ws.Range("H72").WrapText = $false
$ws.Range("H83").Value = 81
$ws.Range("H83").WrapText = $false
$ws.Range("H91").Value = 89
$ws.Range("H91").WrapText = $false
$ws.Range("H99").Value = 97
$ws.Range("H99").WrapText = $false
$ws.Range("H117").Value = 115
$ws.Range("H117").WrapText = $false
$ws.Range("H126").Value = 124
$ws.Range("H126").WrapText = $false
$ws.Range("H135").Value = 133
$ws.Range("H135").WrapText = $false

# Cells that just get a new value, formatting unchanged
$ws.Range("H5").Value = 3
$ws.Range("H7").Value = 5
$ws.Range("H8").Value = 6
$ws.Range("H10").Value = 8
$ws.Range("H13").Value = 11
$ws.Range("H15").Value = 13
$ws.Range("H16").Value = 14
$ws.Range("H20").Value = 18
$ws.Range("H23").Value = 21
$ws.Range("H25").Value = 23
$ws.Range("H26").Value = 24
$ws.Range("H28").Value = 26
$ws.Range("H30").Value = 28
$ws.Range("H31").Value = 29
$ws.Range("H35").Value = 33
$ws.Range("H36").Value = 34
$ws.Range("H38").Value = 36
$ws.Range("H39").Value = 37
$ws.Range("H42").Value = 40
$ws.Range("H44").Value = 42
$ws.Range("H45").Value = 43
$ws.Range("H47").Value = 45
$ws.Range("H48").Value = 46
$ws.Range("H51").Value = 49
$ws.Range("H53").Value = 51
$ws.Range("H54").Value = 52
$ws.Range("H56").Value = 54
$ws.Range("H57").Value = 55
$ws.Range("H60").Value = 58
$ws.Range("H62").Value = 60
$ws.Range("H63").Value = 61
$ws.Range("H65").Value = 63
$ws.Range("H66").Value = 64
$ws.Range("H68").Value = 66
$ws.Range("H69").Value = 67
$ws.Range("H71").Value = 69
$ws.Range("H74").Value = 72
$ws.Range("H77").Value = 75
$ws.Range("H78").Value = 76
$ws.Range("H81").Value = 79
$ws.Range("H82").Value = 80
$ws.Range("H85").Value = 83
$ws.Range("H86").Value = 84
$ws.Range("H87").Value = 85
$ws.Range("H89").Value = 87
$ws.Range("H93").Value = 91
$ws.Range("H94").Value = 92
$ws.Range("H95").Value = 93
$ws.Range("H97").Value = 95
$ws.Range("H101").Value = 99
$ws.Range("H102").Value = 100
$ws.Range("H103").Value = 101
$ws.Range("H105").Value = 103
$ws.Range("H106").Value = 104
$ws.Range("H109").Value = 107
$ws.Range("H110").Value = 108
$ws.Range("H111").Value = 109
$ws.Range("H113").Value = 111
$ws.Range("H114").Value = 112
$ws.Range("H115").Value = 113
$ws.Range("H118").Value = 116
$ws.Range("H119").Value = 117
$ws.Range("H121").Value = 119
$ws.Range("H122").Value = 120
$ws.Range("H123").Value = 121
$ws.Range("H127").Value = 125
$ws.Range("H129").Value = 127
$ws.Range("H130").Value = 128
$ws.Range("H131").Value = 129
$ws.Range("H133").Value = 131
$ws.Range("H137").Value = 135
$ws.Range("H138").Value = 136
$ws.Range("H139").Value = 137
$ws.Range("H141").Value = 139
$ws.Range("H142").Value = 140
$ws.Range("H145").Value = 143
$ws.Range("H146").Value = 144
$ws.Range("H147").Value = 145
$ws.Range("H149").Value = 147
$ws.Range("H150").Value = 148
$ws.Range("H151").Value = 149

# --- View state: selection + window geometry ---
$ws.Range("K23").Select()
